# CUN_1_Control de stock de central.docx
#
# 1) Delete the whole list-item paragraph "Demora la solicitud porque no
#    hay mercadería disponible en el depósito." (including its paragraph
#    mark), which merges it away and leaves the following list item
#    ("Demora de la solicitud por agentes externos...") as the
#    surviving paragraph in that slot.
# 2) Word's auto-tracked "_GoBack" bookmark (marking the last edit
#    position) moves from its old spot (an empty paragraph further down,
#    near the inserted picture) to the start of the paragraph that now
#    sits where the deleted sentence used to be.

$d = $word.ActiveDocument

# --- Step 1: remove the "Demora la solicitud porque no hay mercadería
#     disponible en el depósito." paragraph entirely (text + paragraph
#     mark), merging it with the paragraph that follows it. ---
$targetText = "Demora la solicitud porque no hay mercadería disponible en el depósito."
$paras = $d.Paragraphs
$deleted = $false
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Text -like "*$targetText*") {
        $p.Range.Delete() | Out-Null
        $deleted = $true
        break
    }
}

# --- Step 2: relocate the "_GoBack" bookmark. ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete() | Out-Null
}

$nextText = "Demora de la solicitud por agentes externos que dificultan su confección y/o entrega."
$paras = $d.Paragraphs
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Text -like "*$nextText*") {
        $bmRange = $d.Range($p.Range.Start, $p.Range.Start)
        $d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
        break
    }
}
